$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# The sheet gains a new row 3 ("comment" / "fisrt slide bar" /
# "second slide bar"), pushing the old rows 3-5 down by one: old row 3
# merges into new row 4 together with two brand-new cells, and old rows
# 4 and 5 become new rows 5 and 6.
#
# Work bottom-up so a source cell is never read after being clobbered.
# Formats are carried along with Copy + PasteSpecial(xlPasteFormats);
# that call has to run AFTER the destination's Value is set, because
# assigning .Value resets a cell's number format/quote-prefix flags.
# ---------------------------------------------------------------------

# --- old row 5 (A5 = "EndTestCase") -> new row 6 (no special style) ---
$ws.Range("A6").Value = $ws.Range("A5").Value2

# --- old row 4 (A4="End", B4/C4/D4 blank-but-styled) -> new row 5 ---
$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial($xlPasteFormats)
$ws.Range("C4").Copy()
$ws.Range("C5").PasteSpecial($xlPasteFormats)
$ws.Range("D4").Copy()
$ws.Range("D5").PasteSpecial($xlPasteFormats)
$ws.Range("A5").Value = $ws.Range("A4").Value2
$ws.Range("A4").ClearContents()

# --- old row 3 (B3="amazonitemadd", D3 blank-but-styled) -> new row 4,
#     plus two brand-new cells (C4="Fire TV", D4 reuses old D3 style) ---
$ws.Range("D3").Copy()
$ws.Range("B4").Value = $ws.Range("B3").Value2
$ws.Range("C4").Value = "Fire TV"
$ws.Range("D4").Value = "Fire TV Cube"
$ws.Range("D4").PasteSpecial($xlPasteFormats)

# --- new row 3 content (replaces old row 3), in a new light-blue font ---
$ws.Range("B3").Value = "comment"
$ws.Range("B3").Font.Color = 15773696
$ws.Range("B3").WrapText = $true
$ws.Range("B3").VerticalAlignment = -4160

$ws.Range("C3").Value = "fisrt slide bar"
$ws.Range("D3").Value = "second slide bar"
$ws.Range("E3").Value = ""
$ws.Range("C3:E3").Font.Color = 15773696

$ws.Range("B3").Select()
